$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.808.76"
$ws.Range("D3").Value = "2.290.03"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.68"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.53"
$ws.Range("E6").Value = "  -2.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.511"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -3.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.59"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.68"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  -2.10%  "
$ws.Range("D15").Value = "2.647.04"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("D16").Value = "2.296.18"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").Value = "42.735.85"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.68"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -2.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.71"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "241.01"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.13"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.81"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.02"
$ws.Range("E31").Value = "  -1.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.76"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.81"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.06"
$ws.Range("E36").Value = "  -6.85%  "
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -1.76%  "
$ws.Range("E40").Value = "  -3.78%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.74"
$ws.Range("E42").Value = "  -0.71%  "
$ws.Range("D43").Value = "2.007.44"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("E44").Value = "  -2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.07"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "2.514.25"
$ws.Range("E49").Value = "  -1.09%  "
$ws.Range("B50").Value = "MultiversX"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.08"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.86"
$ws.Range("E51").Value = "  -5.82%  "
